$wb = $excel.ActiveWorkbook

# --- Monsters sheet: add new monster entries and monster-group entries ---
$wsMonsters = $wb.Worksheets.Item("Monsters")

# Row 4 (existing H4/I4 stay, add A4/B4/C4)
$wsMonsters.Range("C4").Value = "Monster in Manyeyes' castle"
$wsMonsters.Range("B4").Value = "Untoter Krieger"
$wsMonsters.Range("A4").Value = 59

# Row 5 (new row entirely)
$wsMonsters.Range("B5").Value = "Untoter Magier"
$wsMonsters.Range("C5").Value = "Monster in Manyeyes' castle"
$wsMonsters.Range("A5").Value = 60
$wsMonsters.Range("H5").Value = 90
$wsMonsters.Range("I6").Value = "2x Untoter Krieger, 1x Untoter Magier"
$wsMonsters.Range("I5").Value = "2x Untoter Krieger"

# Row 6 (new row, only H/I columns)
$wsMonsters.Range("H6").Value = 91

# Row 7 (new row, only H/I columns)
$wsMonsters.Range("H7").Value = 92
$wsMonsters.Range("I7").Value = "3x Untoter Krieger, 2x Untoter Magier"

# Widen column C to fit the new longer description text
$wsMonsters.Columns.Item(3).ColumnWidth = 25.6

# --- View/selection state changes ---
# Previously active sheet (GlobalVars) moves its selection, and loses tabSelected
# once another sheet becomes the active tab.
$wsGlobalVars = $wb.Worksheets.Item("GlobalVars")
$wsGlobalVars.Activate()
$wsGlobalVars.Range("D16").Select() | Out-Null

# Monsters becomes the new active sheet/tab.
$wsMonsters.Activate()
$wsMonsters.Range("I8").Select() | Out-Null
